$d = $word.ActiveDocument
$t = $d.Tables(1)

$cell = $t.Cell(1,1)
$cell.Range.Text = "95÷5=19, 0"

$cell = $t.Cell(1,2)
$cell.Range.Text = "64÷6=10, 4"

$cell = $t.Cell(1,3)
$cell.Range.Text = "65÷4=16, 1"

$cell = $t.Cell(1,4)
$cell.Range.Text = "74÷3=24, 2"

$cell = $t.Cell(1,5)
$cell.Range.Text = "56÷2=28, 0"

$cell = $t.Cell(5,1)
$cell.Range.Text = "38÷2=19, 0"

$cell = $t.Cell(5,2)
$cell.Range.Text = "15÷8=1, 7"

$cell = $t.Cell(5,3)
$cell.Range.Text = "75÷9=8, 3"

$cell = $t.Cell(5,4)
$cell.Range.Text = "49÷7=7, 0"

$cell = $t.Cell(5,5)
$cell.Range.Text = "19÷6=3, 1"

$cell = $t.Cell(9,1)
$cell.Range.Text = "16÷6=2, 4"

$cell = $t.Cell(9,2)
$cell.Range.Text = "64÷2=32, 0"

$cell = $t.Cell(9,3)
$cell.Range.Text = "91÷4=22, 3"

$cell = $t.Cell(9,4)
$cell.Range.Text = "99÷5=19, 4"

$cell = $t.Cell(9,5)
$cell.Range.Text = "89÷5=17, 4"

$cell = $t.Cell(13,1)
$cell.Range.Text = "74÷7=10, 4"

$cell = $t.Cell(13,2)
$cell.Range.Text = "34÷6=5, 4"

$cell = $t.Cell(13,3)
$cell.Range.Text = "63÷2=31, 1"

$cell = $t.Cell(13,4)
$cell.Range.Text = "82÷2=41, 0"

$cell = $t.Cell(13,5)
$cell.Range.Text = "62÷9=6, 8"

$cell = $t.Cell(17,1)
$cell.Range.Text = "35÷8=4, 3"

$cell = $t.Cell(17,2)
$cell.Range.Text = "88÷7=12, 4"

$cell = $t.Cell(17,3)
$cell.Range.Text = "54÷3=18, 0"

$cell = $t.Cell(17,4)
$cell.Range.Text = "55÷9=6, 1"

$cell = $t.Cell(17,5)
$cell.Range.Text = "67÷3=22, 1"
